$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4
$ws.Range("B4").Value = 0.6
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.9
$ws.Range("R4").Value = 0.4285714285714285
$ws.Range("T4").Value = 0
$ws.Range("U4").Value = 0
$ws.Range("V4").Value = 1
$ws.Range("X4").Value = 1
$ws.Range("AG4").Value = 8

# Row 29
$ws.Range("B29").Value = 0.6
$ws.Range("C29").Value = 0
$ws.Range("E29").Value = 1
$ws.Range("F29").Value = 0.8
$ws.Range("G29").Value = 0.7
$ws.Range("I29").Value = 1
$ws.Range("R29").Value = 0.7142857142857143
$ws.Range("S29").Value = 0
$ws.Range("T29").Value = 0
$ws.Range("V29").Value = 0.6
$ws.Range("W29").Value = 0
$ws.Range("X29").Value = 1
$ws.Range("AG29").Value = 8

# Row 42
$ws.Range("B42").Value = 0.6
$ws.Range("C42").Value = 0
$ws.Range("R42").Value = 1
$ws.Range("S42").Value = 0
$ws.Range("T42").Value = 0
$ws.Range("U42").Value = 0
$ws.Range("V42").Value = 1
$ws.Range("W42").Value = 0
$ws.Range("X42").Value = 1
$ws.Range("AG42").Value = 15
$ws.Range("AH42").Value = "MM"
